$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: CreateRecipient
#   - C2 value updated (Fax Number 9987288 -> 918279)
#   - new row 3 added (Data2 / (blank) / 91827)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CreateRecipient")
$ws1.Range("C2").Value = "'918279"
$ws1.Range("A3").Value = "Data2"
$ws1.Range("C3").Value = "'91827"

# ---------------------------------------------------------------------------
# Sheet2: EditRecipient
#   - C2 value updated (Fax Number 9987288 -> 91827)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("EditRecipient")
$ws2.Range("C2").Value = "'91827"
$ws2.Columns.Item(2).ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------------
# Sheet3: DeleteRecipient
#   - C2 value updated (Number 9987288 -> 91827)
#   - new row 5 added (Data2 / (blank) / 91827 / deleted)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DeleteRecipient")
$ws3.Range("C2").Value = "'91827"
$ws3.Range("A5").Value = "Data2"
$ws3.Range("C5").Value = "'91827"
$ws3.Range("D5").Value = "deleted"
$ws3.Columns.Item(4).ColumnWidth = 15.5

# ---------------------------------------------------------------------------
# Sheet4: AddressCreate
#   - B2 (FaxLine) 1000 -> 10384
#   - C2 (Recipient) TrialData Recipient<9987288> -> TrialData Recipient<91827>
#   - D2 (Recipient1) sample2 delete<12345> -> Palak Garg<9917186286>
#   - E2 (Recipients) TrialData Recipient<9987288>,sample2 delete<12345>
#         -> TrialData Recipient<991827>,Palak Garg<9917186286>
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("AddressCreate")
$ws4.Range("B2").Value = "'10384"
$ws4.Range("C2").Value = "TrialData Recipient<91827>"
$ws4.Range("D2").Value = "Palak Garg<9917186286>"
$ws4.Range("E2").Value = "TrialData Recipient<991827>,Palak Garg<9917186286>"
$ws4.Columns.Item(4).ColumnWidth = 27.5

# ---------------------------------------------------------------------------
# Sheet5: EditAddressBook
#   - B2 (FaxLine) 1000 -> 10384
#   - C2 (Recipients) TrialData Recipient<9987288>,sample2 delete<12345>
#         -> TrialData Recipient<91827>,Palak Garg<9917186286>
#   - E2 (Updated Recipients) Fax Address Updated Recipient<9987288>,SampleData<123467>
#         -> Recepient Updated Recipient<91827>,Palak Garg<9917186286>
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("EditAddressBook")
$ws5.Range("B2").Value = "'10384"
$ws5.Range("C2").Value = "TrialData Recipient<91827>,Palak Garg<9917186286>"
$ws5.Range("E2").Value = "Recepient Updated Recipient<91827>,Palak Garg<9917186286>"

# ---------------------------------------------------------------------------
# Sheet6: DeleteAddressBook
#   - B2 (FaxLine) 1000 -> 10384
#   - C2 (Recipients) Fax Address Updated Recipient<9987288>,SampleData<123467>
#         -> Recepient Updated Recipient<91827>,Palak Garg<9917186286>
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("DeleteAddressBook")
$ws6.Range("B2").Value = "'10384"
$ws6.Range("C2").Value = "Recepient Updated Recipient<91827>,Palak Garg<9917186286>"

# ---------------------------------------------------------------------------
# View state: selections per-sheet + active sheet/tab
# ---------------------------------------------------------------------------
$ws2.Range("C8").Select()
$ws3.Range("A7").Select()
$ws4.Range("E6").Select()
$ws5.Range("E6").Select()
$ws6.Range("E5").Select()

$ws1.Activate()
$ws1.Range("F14").Select()
